# This script inserts a new week of price data (2 rows: "Primera" and
# "Segunda" quality records) at the top of the data table on Sheet1,
# pushing the existing records down by two rows. The new rows duplicate
# the (then) first data records except for an updated date (Fecha) and
# updated Volumen values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 253 (shifts rows 253:320 down to 255:322)
$ws.Rows("253:254").Insert()

# --- New row 253 ("Primera" quality) ---
$ws.Range("A253").Value = 1
$ws.Range("B253").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C253").Value = "Arica y Parinacota"
$ws.Range("D253").Value = 44736
$ws.Range("E253").Value = 15
$ws.Range("F253").Value = 100114014
$ws.Range("G253").Value = "Betarraga"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 1000
$ws.Range("K253").Value = 450
$ws.Range("L253").Value = 500
$ws.Range("M253").Value = 475
$ws.Range("N253").Value = "$/paquete 4 unidades"
$ws.Range("O253").Value = "Región de Arica y Parinacota"
$ws.Range("P253").Value = 119
$ws.Range("Q253").Value = 4
$ws.Range("R253").Value = "Hortaliza"

# --- New row 254 ("Segunda" quality) ---
$ws.Range("A254").Value = 1
$ws.Range("B254").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C254").Value = "Arica y Parinacota"
$ws.Range("D254").Value = 44736
$ws.Range("E254").Value = 15
$ws.Range("F254").Value = 100114014
$ws.Range("G254").Value = "Betarraga"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Segunda"
$ws.Range("J254").Value = 1200
$ws.Range("K254").Value = 450
$ws.Range("L254").Value = 500
$ws.Range("M254").Value = 475
$ws.Range("N254").Value = "$/paquete 5 unidades"
$ws.Range("O254").Value = "Región de Arica y Parinacota"
$ws.Range("P254").Value = 95
$ws.Range("Q254").Value = 5
$ws.Range("R254").Value = "Hortaliza"
